$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 ("Graficar_HL" run) label: renamed from 'A' / 'prueba' to the new run name
$ws.Range("A2").Value = "v8-no_fl"
$ws.Range("B2").Value = "Versión 8 SIN FLUORESCENCIA"

# Fixed the table reading: CHL_min/max and NAP_min/max/factor corrected, and
# CHL_factor / NAP_factor now hold plain numbers instead of the step-list text
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 100
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 500
$ws.Range("H2").Value = 1.5

# E2 (CHL_factor) becomes a plain number too; match the formatting used by
# the other "General" numeric cells in the row (e.g. V2) instead of the old
# comment-style format
$ws.Range("E2").Value = 1.5
$ws.Range("V2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# CDOM443_factor keeps the descriptive step-list text, just with the updated values
$ws.Range("K2").Value = "0, 0.1, 1.0, 5.0, 10.0"

# CDOM_steps text, entered as text (not parsed as a number)
$ws.Range("T2").Value = "0.01, 0.018, 0.03"

# Restore the view / selection state used while editing the table
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("T2").Select() | Out-Null
